$wb = $excel.ActiveWorkbook

# Overview sheet: row 7 (bc712ebf-9f83-4070-9206-fc6986c2fb9f.md) -
# "Latest HO Xliff Generate Date" (column G) gets a refreshed timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-22 14:46:45"

# zh-cn sheet: row 7 (bc712ebf-9f83-4070-9206-fc6986c2fb9f.md) -
# "Latest Handoff Datetime" (column H) gets a refreshed timestamp.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-22 14:46:40"

# de-de sheet: row 6 (6e8b221b-b7c4-4b1c-ad67-3a27c471484a.md) and row 7
# (bc712ebf-9f83-4070-9206-fc6986c2fb9f.md) -
# "Latest Handoff Datetime" (column H) gets a refreshed timestamp.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-22 14:46:45"
$wsDeDe.Range("H7").Value = "2016-08-22 14:46:45"
